$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20
$ws.Range("F20").Value = "Winterthur"
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = "Grasshoppers"
$ws.Range("I20").Value = 1
$ws.Range("J20").Value = 2.05
$ws.Range("K20").Value = "06/08/2023 16:42"
$ws.Range("L20").Value = 2.25
$ws.Range("M20").Value = "12/08/2023 17:59"
$ws.Range("N20").Value = 3.64
$ws.Range("O20").Value = "06/08/2023 16:42"
$ws.Range("P20").Value = 3.69
$ws.Range("Q20").Value = "12/08/2023 17:59"
$ws.Range("R20").Value = 3.33
$ws.Range("S20").Value = "06/08/2023 16:42"
$ws.Range("T20").Value = 3.14
$ws.Range("U20").Value = "12/08/2023 17:59"
$ws.Range("V20").Value = "https://www.betexplorer.com/football/switzerland/super-league/winterthur-grasshoppers/QRorOCxF/"

# Row 21
$ws.Range("F21").Value = "Lausanne Ouchy"
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = "Zurich"
$ws.Range("I21").Value = 3
$ws.Range("J21").Value = 3.49
$ws.Range("K21").Value = "08/08/2023 15:12"
$ws.Range("L21").Value = 3.72
$ws.Range("M21").Value = "12/08/2023 17:58"
$ws.Range("N21").Value = 3.74
$ws.Range("O21").Value = "08/08/2023 15:12"
$ws.Range("P21").Value = 3.91
$ws.Range("Q21").Value = "12/08/2023 17:59"
$ws.Range("R21").Value = 1.96
$ws.Range("S21").Value = "08/08/2023 15:12"
$ws.Range("T21").Value = 1.96
$ws.Range("U21").Value = "12/08/2023 17:58"
$ws.Range("V21").Value = "https://www.betexplorer.com/football/switzerland/super-league/lausanne-ouchy-zurich/I1nvPhN8/"

# Row 24
$ws.Range("F24").Value = "Lugano"
$ws.Range("G24").Value = 6
$ws.Range("H24").Value = "Yverdon"
$ws.Range("I24").Value = 1
$ws.Range("J24").Value = 1.36
$ws.Range("K24").Value = "08/08/2023 15:12"
$ws.Range("L24").Value = 1.43
$ws.Range("M24").Value = "13/08/2023 16:21"
$ws.Range("N24").Value = 5.04
$ws.Range("O24").Value = "08/08/2023 15:12"
$ws.Range("P24").Value = 5.15
$ws.Range("Q24").Value = "13/08/2023 16:21"
$ws.Range("R24").Value = 7.46
$ws.Range("S24").Value = "08/08/2023 15:12"
$ws.Range("T24").Value = 6.98
$ws.Range("U24").Value = "13/08/2023 16:21"
$ws.Range("V24").Value = "https://www.betexplorer.com/football/switzerland/super-league/lugano-yverdon/O2RIHAFr/"

# Row 25
$ws.Range("F25").Value = "Luzern"
$ws.Range("G25").Value = 1
$ws.Range("H25").Value = "Young Boys"
$ws.Range("I25").Value = 1
$ws.Range("J25").Value = 3.04
$ws.Range("K25").Value = "06/08/2023 16:42"
$ws.Range("L25").Value = 2.93
$ws.Range("M25").Value = "13/08/2023 16:29"
$ws.Range("N25").Value = 3.81
$ws.Range("O25").Value = "06/08/2023 16:42"
$ws.Range("P25").Value = 3.67
$ws.Range("Q25").Value = "13/08/2023 16:29"
$ws.Range("R25").Value = 2.21
$ws.Range("S25").Value = "06/08/2023 16:42"
$ws.Range("T25").Value = 2.38
$ws.Range("U25").Value = "13/08/2023 16:29"
$ws.Range("V25").Value = "https://www.betexplorer.com/football/switzerland/super-league/luzern-young-boys/IPQMGUUl/"

# Row 26
$ws.Range("F26").Value = "Yverdon"
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = "Servette"
$ws.Range("I26").Value = 1
$ws.Range("J26").Value = 3.88
$ws.Range("K26").Value = "18/08/2023 18:42"
$ws.Range("L26").Value = 4.9
$ws.Range("M26").Value = "26/08/2023 17:57"
$ws.Range("N26").Value = 3.96
$ws.Range("O26").Value = "18/08/2023 18:42"
$ws.Range("P26").Value = 4.46
$ws.Range("Q26").Value = "26/08/2023 17:57"
$ws.Range("R26").Value = 1.81
$ws.Range("S26").Value = "18/08/2023 18:42"
$ws.Range("T26").Value = 1.64
$ws.Range("U26").Value = "26/08/2023 17:45"
$ws.Range("V26").Value = "https://www.betexplorer.com/football/switzerland/super-league/yverdon-servette/dGPQFlpe/"

# Row 27
$ws.Range("F27").Value = "Zurich"
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = "St. Gallen"
$ws.Range("I27").Value = 1
$ws.Range("J27").Value = 2.39
$ws.Range("K27").Value = "18/08/2023 18:42"
$ws.Range("L27").Value = 2.32
$ws.Range("M27").Value = "26/08/2023 17:59"
$ws.Range("N27").Value = 3.58
$ws.Range("O27").Value = "18/08/2023 18:42"
$ws.Range("P27").Value = 3.62
$ws.Range("Q27").Value = "26/08/2023 17:59"
$ws.Range("R27").Value = 2.9
$ws.Range("S27").Value = "18/08/2023 18:42"
$ws.Range("T27").Value = 3.07
$ws.Range("U27").Value = "26/08/2023 17:59"
$ws.Range("V27").Value = "https://www.betexplorer.com/football/switzerland/super-league/zurich-st-gallen/EiKVE8a1/"

# Row 46
$ws.Range("F46").Value = "Basel"
$ws.Range("G46").Value = 1
$ws.Range("H46").Value = "Luzern"
$ws.Range("I46").Value = 1
$ws.Range("J46").Value = 2.34
$ws.Range("K46").Value = "24/09/2023 15:42"
$ws.Range("L46").Value = 2.57
$ws.Range("M46").Value = "28/09/2023 20:03"
$ws.Range("N46").Value = 3.65
$ws.Range("O46").Value = "24/09/2023 15:42"
$ws.Range("P46").Value = 3.82
$ws.Range("Q46").Value = "28/09/2023 20:03"
$ws.Range("R46").Value = 2.79
$ws.Range("S46").Value = "24/09/2023 15:42"
$ws.Range("T46").Value = 2.62
$ws.Range("U46").Value = "28/09/2023 20:03"
$ws.Range("V46").Value = "https://www.betexplorer.com/football/switzerland/super-league/basel-luzern/YVDhQd4U/"

# Row 47
$ws.Range("F47").Value = "Lausanne Ouchy"
$ws.Range("G47").Value = 1
$ws.Range("H47").Value = "Yverdon"
$ws.Range("I47").Value = 1
$ws.Range("J47").Value = 1.88
$ws.Range("K47").Value = "24/09/2023 22:12"
$ws.Range("L47").Value = 2.23
$ws.Range("M47").Value = "28/09/2023 20:29"
$ws.Range("N47").Value = 3.89
$ws.Range("O47").Value = "24/09/2023 22:12"
$ws.Range("P47").Value = 3.76
$ws.Range("Q47").Value = "28/09/2023 20:28"
$ws.Range("R47").Value = 3.66
$ws.Range("S47").Value = "24/09/2023 22:12"
$ws.Range("T47").Value = 3.13
$ws.Range("U47").Value = "28/09/2023 20:29"
$ws.Range("V47").Value = "https://www.betexplorer.com/football/switzerland/super-league/lausanne-ouchy-yverdon/Aikbotzt/"

# Row 52
$ws.Range("F52").Value = "Basel"
$ws.Range("G52").Value = 0
$ws.Range("H52").Value = "Lausanne Ouchy"
$ws.Range("I52").Value = 3
$ws.Range("J52").Value = 1.66
$ws.Range("K52").Value = "28/09/2023 19:42"
$ws.Range("L52").Value = 1.78
$ws.Range("M52").Value = "01/10/2023 16:29"
$ws.Range("N52").Value = 4.21
$ws.Range("O52").Value = "28/09/2023 19:42"
$ws.Range("P52").Value = 4.12
$ws.Range("Q52").Value = "01/10/2023 16:29"
$ws.Range("R52").Value = 4.81
$ws.Range("S52").Value = "28/09/2023 19:42"
$ws.Range("T52").Value = 4.3
$ws.Range("U52").Value = "01/10/2023 16:26"
$ws.Range("V52").Value = "https://www.betexplorer.com/football/switzerland/super-league/basel-lausanne-ouchy/hzoJtIlB/"

# Row 53
$ws.Range("F53").Value = "Luzern"
$ws.Range("G53").Value = 1
$ws.Range("H53").Value = "Zurich"
$ws.Range("I53").Value = 4
$ws.Range("J53").Value = 2.26
$ws.Range("K53").Value = "28/09/2023 19:42"
$ws.Range("L53").Value = 2.42
$ws.Range("M53").Value = "01/10/2023 16:20"
$ws.Range("N53").Value = 3.65
$ws.Range("O53").Value = "28/09/2023 19:42"
$ws.Range("P53").Value = 3.61
$ws.Range("Q53").Value = "01/10/2023 16:20"
$ws.Range("R53").Value = 2.9
$ws.Range("S53").Value = "28/09/2023 19:42"
$ws.Range("T53").Value = 2.92
$ws.Range("U53").Value = "01/10/2023 16:20"
$ws.Range("V53").Value = "https://www.betexplorer.com/football/switzerland/super-league/luzern-zurich/zazOux4H/"

# Row 58
$ws.Range("F58").Value = "Lausanne"
$ws.Range("G58").Value = 3
$ws.Range("H58").Value = "Luzern"
$ws.Range("I58").Value = 1
$ws.Range("J58").Value = 2.33
$ws.Range("K58").Value = "03/10/2023 02:13"
$ws.Range("L58").Value = 2.67
$ws.Range("M58").Value = "08/10/2023 16:29"
$ws.Range("N58").Value = 3.63
$ws.Range("O58").Value = "03/10/2023 02:13"
$ws.Range("P58").Value = 3.54
$ws.Range("Q58").Value = "08/10/2023 16:24"
$ws.Range("R58").Value = 2.81
$ws.Range("S58").Value = "03/10/2023 02:13"
$ws.Range("T58").Value = 2.66
$ws.Range("U58").Value = "08/10/2023 16:29"
$ws.Range("V58").Value = "https://www.betexplorer.com/football/switzerland/super-league/lausanne-luzern/QsIqcwZi/"

# Row 59
$ws.Range("F59").Value = "Young Boys"
$ws.Range("G59").Value = 3
$ws.Range("H59").Value = "Basel"
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 1.47
$ws.Range("K59").Value = "01/10/2023 15:42"
$ws.Range("L59").Value = 1.51
$ws.Range("M59").Value = "08/10/2023 16:19"
$ws.Range("N59").Value = 4.85
$ws.Range("O59").Value = "01/10/2023 15:42"
$ws.Range("P59").Value = 4.89
$ws.Range("Q59").Value = "08/10/2023 16:29"
$ws.Range("R59").Value = 5.56
$ws.Range("S59").Value = "01/10/2023 15:42"
$ws.Range("T59").Value = 5.82
$ws.Range("U59").Value = "08/10/2023 16:29"
$ws.Range("V59").Value = "https://www.betexplorer.com/football/switzerland/super-league/young-boys-basel/K4Hmdclb/"

# Row 60
$ws.Range("F60").Value = "Young Boys"
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = "Zurich"
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 1.76
$ws.Range("K60").Value = "11/10/2023 13:43"
$ws.Range("L60").Value = 2.03
$ws.Range("M60").Value = "21/10/2023 17:59"
$ws.Range("N60").Value = 4.21
$ws.Range("O60").Value = "11/10/2023 13:43"
$ws.Range("P60").Value = 3.7
$ws.Range("Q60").Value = "21/10/2023 17:59"
$ws.Range("R60").Value = 3.84
$ws.Range("S60").Value = "11/10/2023 13:43"
$ws.Range("T60").Value = 3.69
$ws.Range("U60").Value = "21/10/2023 17:59"
$ws.Range("V60").Value = "https://www.betexplorer.com/football/switzerland/super-league/young-boys-zurich/GWAdfyJA/"

# Row 61
$ws.Range("F61").Value = "Lausanne Ouchy"
$ws.Range("G61").Value = 2
$ws.Range("H61").Value = "Lausanne"
$ws.Range("I61").Value = 2
$ws.Range("J61").Value = 2.88
$ws.Range("K61").Value = "11/10/2023 13:43"
$ws.Range("L61").Value = 2.88
$ws.Range("M61").Value = "21/10/2023 17:55"
$ws.Range("N61").Value = 3.47
$ws.Range("O61").Value = "11/10/2023 13:43"
$ws.Range("P61").Value = 3.48
$ws.Range("Q61").Value = "21/10/2023 17:55"
$ws.Range("R61").Value = 2.35
$ws.Range("S61").Value = "11/10/2023 13:43"
$ws.Range("T61").Value = 2.51
$ws.Range("U61").Value = "21/10/2023 17:55"
$ws.Range("V61").Value = "https://www.betexplorer.com/football/switzerland/super-league/lausanne-ouchy-lausanne/beGieH34/"

# Row 66
$ws.Range("F66").Value = "St. Gallen"
$ws.Range("G66").Value = 3
$ws.Range("H66").Value = "Grasshoppers"
$ws.Range("I66").Value = 1
$ws.Range("J66").Value = 1.39
$ws.Range("K66").Value = "22/10/2023 16:42"
$ws.Range("L66").Value = 1.62
$ws.Range("M66").Value = "28/10/2023 17:36"
$ws.Range("N66").Value = 5.14
$ws.Range("O66").Value = "22/10/2023 16:42"
$ws.Range("P66").Value = 4.56
$ws.Range("Q66").Value = "28/10/2023 17:59"
$ws.Range("R66").Value = 6.33
$ws.Range("S66").Value = "22/10/2023 16:42"
$ws.Range("T66").Value = 5.02
$ws.Range("U66").Value = "28/10/2023 17:59"
$ws.Range("V66").Value = "https://www.betexplorer.com/football/switzerland/super-league/st-gallen-grasshoppers/fNVV8zY3/"

# Row 67
$ws.Range("F67").Value = "Yverdon"
$ws.Range("G67").Value = 1
$ws.Range("H67").Value = "Winterthur"
$ws.Range("I67").Value = 1
$ws.Range("J67").Value = 2.65
$ws.Range("K67").Value = "22/10/2023 16:42"
$ws.Range("L67").Value = 2.47
$ws.Range("M67").Value = "28/10/2023 17:59"
$ws.Range("N67").Value = 3.64
$ws.Range("O67").Value = "22/10/2023 16:42"
$ws.Range("P67").Value = 3.73
$ws.Range("Q67").Value = "28/10/2023 17:57"
$ws.Range("R67").Value = 2.44
$ws.Range("S67").Value = "22/10/2023 16:42"
$ws.Range("T67").Value = 2.78
$ws.Range("U67").Value = "28/10/2023 17:59"
$ws.Range("V67").Value = "https://www.betexplorer.com/football/switzerland/super-league/yverdon-winterthur/8QZZ7fmA/"

# Row 70
$ws.Range("F70").Value = "Lugano"
$ws.Range("G70").Value = 1
$ws.Range("H70").Value = "Young Boys"
$ws.Range("I70").Value = 1
$ws.Range("J70").Value = 2.95
$ws.Range("K70").Value = "22/10/2023 20:15"
$ws.Range("L70").Value = 2.73
$ws.Range("M70").Value = "29/10/2023 16:28"
$ws.Range("N70").Value = 3.68
$ws.Range("O70").Value = "22/10/2023 20:15"
$ws.Range("P70").Value = 3.66
$ws.Range("Q70").Value = "29/10/2023 16:28"
$ws.Range("R70").Value = 2.22
$ws.Range("S70").Value = "22/10/2023 20:15"
$ws.Range("T70").Value = 2.53
$ws.Range("U70").Value = "29/10/2023 16:29"
$ws.Range("V70").Value = "https://www.betexplorer.com/football/switzerland/super-league/lugano-young-boys/Y7sNAd3i/"

# Row 71
$ws.Range("F71").Value = "Servette"
$ws.Range("G71").Value = 4
$ws.Range("H71").Value = "Luzern"
$ws.Range("I71").Value = 2
$ws.Range("J71").Value = 1.99
$ws.Range("K71").Value = "22/10/2023 20:15"
$ws.Range("L71").Value = 1.93
$ws.Range("M71").Value = "29/10/2023 16:21"
$ws.Range("N71").Value = 3.86
$ws.Range("O71").Value = "22/10/2023 20:15"
$ws.Range("P71").Value = 3.87
$ws.Range("Q71").Value = "29/10/2023 16:21"
$ws.Range("R71").Value = 3.52
$ws.Range("S71").Value = "22/10/2023 20:15"
$ws.Range("T71").Value = 3.9
$ws.Range("U71").Value = "29/10/2023 16:21"
$ws.Range("V71").Value = "https://www.betexplorer.com/football/switzerland/super-league/servette-luzern/CbtR9GIc/"

# Row 72
$ws.Range("F72").Value = "Winterthur"
$ws.Range("G72").Value = 1
$ws.Range("H72").Value = "Young Boys"
$ws.Range("I72").Value = 4
$ws.Range("J72").Value = 3.25
$ws.Range("K72").Value = "29/10/2023 16:42"
$ws.Range("L72").Value = 3.47
$ws.Range("M72").Value = "04/11/2023 17:57"
$ws.Range("N72").Value = 4.23
$ws.Range("O72").Value = "29/10/2023 16:42"
$ws.Range("P72").Value = 3.89
$ws.Range("Q72").Value = "04/11/2023 17:57"
$ws.Range("R72").Value = 1.93
$ws.Range("S72").Value = "29/10/2023 16:42"
$ws.Range("T72").Value = 2.04
$ws.Range("U72").Value = "04/11/2023 17:57"
$ws.Range("V72").Value = "https://www.betexplorer.com/football/switzerland/super-league/winterthur-young-boys/WOvFbjAd/"

# Row 73
$ws.Range("F73").Value = "Lausanne"
$ws.Range("G73").Value = 3
$ws.Range("H73").Value = "Lugano"
$ws.Range("I73").Value = 1
$ws.Range("J73").Value = 2.24
$ws.Range("K73").Value = "29/10/2023 16:42"
$ws.Range("L73").Value = 2.34
$ws.Range("M73").Value = "04/11/2023 17:52"
$ws.Range("N73").Value = 3.67
$ws.Range("O73").Value = "29/10/2023 16:42"
$ws.Range("P73").Value = 3.67
$ws.Range("Q73").Value = "04/11/2023 17:52"
$ws.Range("R73").Value = 3.08
$ws.Range("S73").Value = "29/10/2023 16:42"
$ws.Range("T73").Value = 3
$ws.Range("U73").Value = "04/11/2023 17:52"
$ws.Range("V73").Value = "https://www.betexplorer.com/football/switzerland/super-league/lausanne-lugano/Mkku5hXS/"

# Row 78
$ws.Range("F78").Value = "St. Gallen"
$ws.Range("G78").Value = 4
$ws.Range("H78").Value = "Winterthur"
$ws.Range("I78").Value = 2
$ws.Range("J78").Value = 1.47
$ws.Range("K78").Value = "05/11/2023 16:42"
$ws.Range("L78").Value = 1.53
$ws.Range("M78").Value = "11/11/2023 17:52"
$ws.Range("N78").Value = 4.98
$ws.Range("O78").Value = "05/11/2023 16:42"
$ws.Range("P78").Value = 4.93
$ws.Range("Q78").Value = "11/11/2023 17:55"
$ws.Range("R78").Value = 5.94
$ws.Range("S78").Value = "05/11/2023 16:42"
$ws.Range("T78").Value = 5.57
$ws.Range("U78").Value = "11/11/2023 17:55"
$ws.Range("V78").Value = "https://www.betexplorer.com/football/switzerland/super-league/st-gallen-winterthur/0rZw351r/"

# Row 79
$ws.Range("F79").Value = "Yverdon"
$ws.Range("G79").Value = 2
$ws.Range("H79").Value = "Lausanne"
$ws.Range("I79").Value = 2
$ws.Range("J79").Value = 3.38
$ws.Range("K79").Value = "05/11/2023 16:42"
$ws.Range("L79").Value = 3.74
$ws.Range("M79").Value = "11/11/2023 17:59"
$ws.Range("N79").Value = 3.81
$ws.Range("O79").Value = "05/11/2023 16:42"
$ws.Range("P79").Value = 3.87
$ws.Range("Q79").Value = "11/11/2023 17:59"
$ws.Range("R79").Value = 2.06
$ws.Range("S79").Value = "05/11/2023 16:42"
$ws.Range("T79").Value = 1.96
$ws.Range("U79").Value = "11/11/2023 17:53"
$ws.Range("V79").Value = "https://www.betexplorer.com/football/switzerland/super-league/yverdon-lausanne/hjio1qWf/"
# New rows 82 and 83: copy formatting from row 81, then set values
$ws.Range("A81:V81").Copy()
$ws.Range("A82").PasteSpecial(-4122)
$ws.Range("A81:V81").Copy()
$ws.Range("A83").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 82
$ws.Range("A82").Value = 81
$ws.Range("B82").Value = "switzerland"
$ws.Range("C82").Value = "super-league"
$ws.Range("D82").Value = "2023-2024"
$ws.Range("E82").Value = 45242.6875
$ws.Range("F82").Value = "Lugano"
$ws.Range("G82").Value = 0
$ws.Range("H82").Value = "Zurich"
$ws.Range("I82").Value = 3
$ws.Range("J82").Value = 2.49
$ws.Range("K82").Value = "05/11/2023 16:42"
$ws.Range("L82").Value = 2.93
$ws.Range("M82").Value = "12/11/2023 16:23"
$ws.Range("N82").Value = 3.51
$ws.Range("O82").Value = "05/11/2023 16:42"
$ws.Range("P82").Value = 3.26
$ws.Range("Q82").Value = "12/11/2023 16:29"
$ws.Range("R82").Value = 2.81
$ws.Range("S82").Value = "05/11/2023 16:42"
$ws.Range("T82").Value = 2.59
$ws.Range("U82").Value = "12/11/2023 16:23"
$ws.Range("V82").Value = "https://www.betexplorer.com/football/switzerland/super-league/lugano-zurich/tfM5LDH9/"

# Row 83
$ws.Range("A83").Value = 82
$ws.Range("B83").Value = "switzerland"
$ws.Range("C83").Value = "super-league"
$ws.Range("D83").Value = "2023-2024"
$ws.Range("E83").Value = 45242.6875
$ws.Range("F83").Value = "Servette"
$ws.Range("G83").Value = 4
$ws.Range("H83").Value = "Basel"
$ws.Range("I83").Value = 1
$ws.Range("J83").Value = 1.67
$ws.Range("K83").Value = "05/11/2023 16:42"
$ws.Range("L83").Value = 1.62
$ws.Range("M83").Value = "12/11/2023 16:27"
$ws.Range("N83").Value = 4.26
$ws.Range("O83").Value = "05/11/2023 16:42"
$ws.Range("P83").Value = 4.37
$ws.Range("Q83").Value = "12/11/2023 16:27"
$ws.Range("R83").Value = 4.72
$ws.Range("S83").Value = "05/11/2023 16:42"
$ws.Range("T83").Value = 5.27
$ws.Range("U83").Value = "12/11/2023 16:25"
$ws.Range("V83").Value = "https://www.betexplorer.com/football/switzerland/super-league/servette-basel/KpL9KXWF/"
